$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_30.2")

# Insert a new row above row 6 (shifts existing rows 6-95 down to 7-96)
$ws.Rows.Item(6).Insert()

# The newly inserted row 6 should pick up the same look as row 8 (which is
# the old row 7, now shifted down one, untouched) -- that's the repeating
# style this table uses for the second row of each 3-row block.
$ws.Range("B8:F8").Copy()
$ws.Range("B6:F6").PasteSpecial(-4122)

# Fill the new row 6 with April 2025 data
$ws.Range("B6").Value = 2025
$ws.Range("C6").Value = "Abr."
$ws.Range("D6").Formula = "=SUM(E6:F6)"
$ws.Range("E6").Value = 5436631
$ws.Range("F6").Value = 5061112

# Update the "Actualización" note text (row 94 in the new layout)
$found = $ws.Cells.Find("Actualización: Marzo 2025.")
if ($found -ne $null) {
    $found.Value = "Actualización: Abril 2025."
}
